# Update column G ("K" - strikeouts) values for rows 2-31 on Sheet1.
# The data was regenerated from source (commit: "regen save_data to use K
# instead of Strike#, regen std/mean, calc and write s_vals"); only the
# numeric values in column G change here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 2
    3  = 7
    4  = 4
    5  = 3
    6  = 3
    7  = 9
    8  = 4
    9  = 6
    10 = 8
    11 = 3
    12 = 4
    13 = 7
    14 = 5
    15 = 8
    16 = 3
    17 = 6
    18 = 8
    19 = 4
    20 = 2
    21 = 3
    22 = 5
    23 = 2
    24 = 4
    25 = 5
    26 = 5
    27 = 6
    28 = 2
    29 = 2
    30 = 6
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
